$wb = $excel.ActiveWorkbook

# ---------- Sheet "Product" ----------
$ws1 = $wb.Worksheets.Item("Product")

# Back up the two header styles currently in use before we start overwriting them:
#   Z1 <- style currently on E2 (bordered + centered "label" look)
#   Z2 <- style currently on A1 (plain bold look)
$ws1.Range("E2").Copy()
$ws1.Range("Z1").PasteSpecial(-4122)
$ws1.Range("A1").Copy()
$ws1.Range("Z2").PasteSpecial(-4122)

# Apply the bordered/centered look (backed up in Z1) to the section header cells
$ws1.Range("Z1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)
$ws1.Range("Z1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)
$ws1.Range("Z1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("Z1").Copy()
$ws1.Range("B2").PasteSpecial(-4122)
$ws1.Range("Z1").Copy()
$ws1.Range("A4").PasteSpecial(-4122)
$ws1.Range("Z1").Copy()
$ws1.Range("A10").PasteSpecial(-4122)

# Apply the plain-bold look (backed up in Z2) to the scale label cell
$ws1.Range("Z2").Copy()
$ws1.Range("E2").PasteSpecial(-4122)

# Remove the scratch cells
$ws1.Range("Z1:Z2").Clear()

# New "Packaging" requirement notes + widen the Notes column (target stored width 16.57)
$ws1.Columns.Item(3).ColumnWidth = 15.74

$ws1.Range("A13").Copy()
$ws1.Range("C13").PasteSpecial(-4122)
$ws1.Range("C13").Value = "Looks sharp"

$ws1.Range("A14").Copy()
$ws1.Range("C14").PasteSpecial(-4122)
$ws1.Range("C14").Value = "Not safe to recycle"

# ---------- Sheet "Company" ----------
$ws2 = $wb.Worksheets.Item("Company")

$ws2.Range("E2").Copy()
$ws2.Range("Z1").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("Z2").PasteSpecial(-4122)

$ws2.Range("Z1").Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Range("Z1").Copy()
$ws2.Range("A10").PasteSpecial(-4122)

$ws2.Range("Z2").Copy()
$ws2.Range("E2").PasteSpecial(-4122)

$ws2.Range("Z1:Z2").Clear()
